$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.064266334283068
$ws.Range("D2").Value = 1.062798811363362
$ws.Range("E2").Value = 1.067969375299551
$ws.Range("F2").Value = 1.072318345708968
$ws.Range("I2").Value = 1.044349348132172
$ws.Range("J2").Value = 1.069227456614976
$ws.Range("K2").Value = 1.065519181323221
$ws.Range("L2").Value = 1.070675813482637
$ws.Range("M2").Value = 1.075013178320608
$ws.Range("N2").Value = 1.070745882245473

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.066141018215524
$ws.Range("D3").Value = 1.064212215937001
$ws.Range("E3").Value = 1.069796595076647
$ws.Range("F3").Value = 1.073885558473027
$ws.Range("I3").Value = 1.044804074581916
$ws.Range("J3").Value = 1.070752902397079
$ws.Range("K3").Value = 1.06674596883758
$ws.Range("L3").Value = 1.072316408056803
$ws.Range("M3").Value = 1.07639526338018
$ws.Range("N3").Value = 1.072273494335557

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.067350788483526
$ws.Range("D4").Value = 1.065123726779319
$ws.Range("E4").Value = 1.070976023199897
$ws.Range("F4").Value = 1.074896778972174
$ws.Range("I4").Value = 1.04509553987373
$ws.Range("J4").Value = 1.071736426952033
$ws.Range("K4").Value = 1.067536205030597
$ws.Range("L4").Value = 1.073374640260272
$ws.Range("M4").Value = 1.077286200280109
$ws.Range("N4").Value = 1.073258415608192

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.067858609720217
$ws.Range("D5").Value = 1.065506206727062
$ws.Range("E5").Value = 1.071471175916684
$ws.Range("F5").Value = 1.075321222334025
$ws.Range("I5").Value = 1.045217412967399
$ws.Range("J5").Value = 1.072149067597703
$ws.Range("K5").Value = 1.067867575673414
$ws.Range("L5").Value = 1.073818736390159
$ws.Range("M5").Value = 1.077659957382229
$ws.Range("N5").Value = 1.073671642250908

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.067943830730767
$ws.Range("D6").Value = 1.065570384979411
$ws.Range("E6").Value = 1.07155427490859
$ws.Range("F6").Value = 1.075392449125881
$ws.Range("I6").Value = 1.045237837500603
$ws.Range("J6").Value = 1.072218303347691
$ws.Range("K6").Value = 1.067923165071477
$ws.Range("L6").Value = 1.073893256515734
$ws.Range("M6").Value = 1.077722666722726
$ws.Range("N6").Value = 1.073740976323603

$ws.Range("B7").Value = 1.019999999999999
$ws.Range("C7").Value = 1.06735757700179
$ws.Range("D7").Value = 1.06512884030398
$ws.Range("E7").Value = 1.070982642099475
$ws.Range("F7").Value = 1.074902453033884
$ws.Range("I7").Value = 1.045097170930214
$ws.Range("J7").Value = 1.07174194393
$ws.Range("K7").Value = 1.067540636124137
$ws.Range("L7").Value = 1.073380577356703
$ws.Range("M7").Value = 1.077291197538451
$ws.Range("N7").Value = 1.073263940420899

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.064900581611754
$ws.Range("D8").Value = 1.063277118473475
$ws.Range("E8").Value = 1.068587504402887
$ws.Range("F8").Value = 1.072848595261348
$ws.Range("I8").Value = 1.044503602270964
$ws.Range("J8").Value = 1.069743730157399
$ws.Range("K8").Value = 1.065934527745849
$ws.Range("L8").Value = 1.07123096145843
$ws.Range("M8").Value = 1.075480964776879
$ws.Range("N8").Value = 1.071262888955544

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.060545142889151
$ws.Range("D9").Value = 1.059990174943042
$ws.Range("E9").Value = 1.064343963599601
$ws.Range("F9").Value = 1.069206834632276
$ws.Range("I9").Value = 1.043436194655049
$ws.Range("J9").Value = 1.066194843928174
$ws.Range("K9").Value = 1.06307642968555
$ws.Range("L9").Value = 1.067416769894175
$ws.Range("M9").Value = 1.07226476122519
$ws.Range("N9").Value = 1.06770896290081

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.057622963506295
$ws.Range("D10").Value = 1.057781966783147
$ws.Range("E10").Value = 1.061498415010269
$ws.Range("F10").Value = 1.066762936738059
$ws.Range("I10").Value = 1.042709841406134
$ws.Range("J10").Value = 1.063809316029272
$ws.Range("K10").Value = 1.061151504477627
$ws.Range("L10").Value = 1.064855317919721
$ws.Range("M10").Value = 1.070102107379366
$ws.Range("N10").Value = 1.065320047278667

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.056352971562898
$ws.Range("D11").Value = 1.056821598659925
$ws.Range("E11").Value = 1.060262105886974
$ws.Range("F11").Value = 1.065700700058196
$ws.Range("I11").Value = 1.042391753031163
$ws.Range("J11").Value = 1.062771498569558
$ws.Range("K11").Value = 1.0603131913717
$ws.Range("L11").Value = 1.063741536213217
$ws.Range("M11").Value = 1.06916108108315
$ws.Range("N11").Value = 1.064280755999122

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.055880516430759
$ws.Range("D12").Value = 1.056464228787609
$ws.Range("E12").Value = 1.0598022389733
$ws.Range("F12").Value = 1.065305518411104
$ws.Range("I12").Value = 1.042273057870849
$ws.Range("J12").Value = 1.062385257547692
$ws.Range("K12").Value = 1.060001067661179
$ws.Range("L12").Value = 1.063327109547507
$ws.Range("M12").Value = 1.068810837764133
$ws.Range("N12").Value = 1.063893966470702

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.055981892761776
$ws.Range("D13").Value = 1.056540915339446
$ws.Range("E13").Value = 1.05990091158171
$ws.Range("F13").Value = 1.065390314527839
$ws.Range("I13").Value = 1.042298543057278
$ws.Range("J13").Value = 1.062468141762116
$ws.Range("K13").Value = 1.060068052838824
$ws.Range("L13").Value = 1.063416038297872
$ws.Range("M13").Value = 1.068885998268484
$ws.Range("N13").Value = 1.063976968390217

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.056313933163986
$ws.Range("D14").Value = 1.056792071631214
$ws.Range("E14").Value = 1.06022410645554
$ws.Range("F14").Value = 1.065668046958411
$ws.Range("I14").Value = 1.042381952765034
$ws.Range("J14").Value = 1.062739587159003
$ws.Range("K14").Value = 1.060287406245919
$ws.Range("L14").Value = 1.063707294344827
$ws.Range("M14").Value = 1.069132144319372
$ws.Range("N14").Value = 1.064248799270703

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.056518417667075
$ws.Range("D15").Value = 1.056946731192326
$ws.Range("E15").Value = 1.06042315117809
$ws.Range("F15").Value = 1.065839084460318
$ws.Range("I15").Value = 1.042433272101658
$ws.Range("J15").Value = 1.06290673377979
$ws.Range("K15").Value = 1.060422459034175
$ws.Range("L15").Value = 1.063886651020152
$ws.Range("M15").Value = 1.06928370923095
$ws.Range("N15").Value = 1.064416183258864

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.057707148437323
$ws.Range("D16").Value = 1.057845613459213
$ws.Range("E16").Value = 1.061580375220822
$ws.Range("F16").Value = 1.066833347916705
$ws.Range("I16").Value = 1.042730876083108
$ws.Range("J16").Value = 1.063878088456579
$ws.Range("K16").Value = 1.061207038022613
$ws.Range("L16").Value = 1.064929136226507
$ws.Range("M16").Value = 1.070164462301303
$ws.Range("N16").Value = 1.06538891737071

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.058451542296113
$ws.Range("D17").Value = 1.058408323389795
$ws.Range("E17").Value = 1.062305141061618
$ws.Range("F17").Value = 1.067455937007
$ws.Range("I17").Value = 1.042916594333395
$ws.Range("J17").Value = 1.064486077588913
$ws.Range("K17").Value = 1.061697886276516
$ws.Range("L17").Value = 1.065581799729565
$ws.Range("M17").Value = 1.070715697420514
$ws.Range("N17").Value = 1.065997769917349

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.058885285657598
$ws.Range("D18").Value = 1.058736138444488
$ws.Range("E18").Value = 1.062727483978536
$ws.Range("F18").Value = 1.067818696660351
$ws.Range("I18").Value = 1.043024576365177
$ws.Range("J18").Value = 1.064840238887722
$ws.Range("K18").Value = 1.061983726871049
$ws.Range("L18").Value = 1.065962039359529
$ws.Range("M18").Value = 1.071036782427267
$ws.Range("N18").Value = 1.066352434165825

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.059033105444852
$ws.Range("D19").Value = 1.058847846824819
$ws.Range("E19").Value = 1.062871424643356
$ws.Range("F19").Value = 1.067942323354602
$ws.Range("I19").Value = 1.043061337238552
$ws.Range("J19").Value = 1.064960919818656
$ws.Range("K19").Value = 1.062081113051591
$ws.Range("L19").Value = 1.066091615839116
$ws.Range("M19").Value = 1.071146189769396
$ws.Range("N19").Value = 1.066473286477523

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.05837172247975
$ws.Range("D20").Value = 1.058347991852108
$ws.Range("E20").Value = 1.062227422143809
$ws.Range("F20").Value = 1.067389179139599
$ws.Range("I20").Value = 1.042896704178337
$ws.Range("J20").Value = 1.064420894695606
$ws.Range("K20").Value = 1.061645270888001
$ws.Range("L20").Value = 1.065511821583105
$ws.Range("M20").Value = 1.070656600848761
$ws.Range("N20").Value = 1.065932494456857

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.056216175694684
$ws.Range("D21").Value = 1.056718130332012
$ws.Range("E21").Value = 1.0601289516374
$ws.Range("F21").Value = 1.065586278921002
$ws.Range("I21").Value = 1.042357405724224
$ws.Range("J21").Value = 1.062659674088447
$ws.Range("K21").Value = 1.060222832594758
$ws.Range("L21").Value = 1.063621546687776
$ws.Range("M21").Value = 1.069059680010516
$ws.Range("N21").Value = 1.064168772714421

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.054856700324547
$ws.Range("D22").Value = 1.055689624077213
$ws.Range("E22").Value = 1.058805808238392
$ws.Range("F22").Value = 1.064449130440029
$ws.Range("I22").Value = 1.042015182751028
$ws.Range("J22").Value = 1.061547979850011
$ws.Range("K22").Value = 1.059324218366196
$ws.Range("L22").Value = 1.062428888938955
$ws.Range("M22").Value = 1.06805154997145
$ws.Range("N22").Value = 1.063055499742653

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.055577789274528
$ws.Range("D23").Value = 1.056235215152601
$ws.Range("E23").Value = 1.059507594114907
$ws.Range("F23").Value = 1.065052300648632
$ws.Range("I23").Value = 1.042196901799414
$ws.Range("J23").Value = 1.062137727926709
$ws.Range("K23").Value = 1.059801000511497
$ws.Range("L23").Value = 1.063061540763758
$ws.Range("M23").Value = 1.06858637105697
$ws.Range("N23").Value = 1.063646085329266

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.058407791005825
$ws.Range("D24").Value = 1.058375254325545
$ws.Range("E24").Value = 1.062262541215262
$ws.Range("F24").Value = 1.067419345338876
$ws.Range("I24").Value = 1.042905692746963
$ws.Range("J24").Value = 1.064450349486547
$ws.Range("K24").Value = 1.061669046947314
$ws.Range("L24").Value = 1.065543443079092
$ws.Range("M24").Value = 1.070683305408508
$ws.Range("N24").Value = 1.06596199107698

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.061674313968651
$ws.Range("D25").Value = 1.060842850752039
$ws.Range("E25").Value = 1.065443854453747
$ws.Range("F25").Value = 1.070151083467423
$ws.Range("I25").Value = 1.043714721079308
$ws.Range("J25").Value = 1.067115702733835
$ws.Range("K25").Value = 1.063818702829042
$ws.Range("L25").Value = 1.068406045880015
$ws.Range("M25").Value = 1.073099431115885
$ws.Range("N25").Value = 1.068631129431599
